# export/excel: export special fields
#
# The "requirements" sheet is renamed to "Requirements" (capitalised tab
# name), and every hyperlink that points back into this sheet (the
# PARENT column's links to the UID column) is repointed so its
# sheet-qualified location string keeps working after the rename.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldName = $ws.Name
$newName = "Requirements"

$ws.Name = $newName

foreach ($hl in $ws.Hyperlinks) {
    $hl.SubAddress = $hl.SubAddress -replace [regex]::Escape("'" + $oldName + "'"), ("'" + $newName + "'")
}
